$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SwateTemplateMetadata")

# Harmonize the "Computational analyses" tag name -> "Computational Analysis"
$ws.Range("B12").Value = "Computational Analysis"

# Add term accession number / term source REF for the "Proteomics" tag (column C)
$ws.Range("C13").Value = "http://purl.obolibrary.org/obo/NCIT_C20085"
$ws.Range("C14").Value = "NCIT"

# Add term accession number / term source REF for the "PRIDE" ER
$ws.Range("B9").Value = "http://purl.obolibrary.org/obo/DPBO_1000098"
$ws.Range("B10").Value = "DPBO"

# Update the active sheet / selection to reflect the edited sheet
$ws.Range("B13").Select()
$ws.Activate()
